$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 16,20
$data[0,0] = "Sending cluster"
$data[0,1] = "Ligand symbol"
$data[0,2] = "Receptor symbol"
$data[0,3] = "Target cluster"
$data[0,4] = "Ligand-expressing cells"
$data[0,5] = "Ligand detection rate"
$data[0,6] = "Ligand average expression value"
$data[0,7] = "Ligand total expression value"
$data[0,8] = "Ligand derived specificity of average expression value"
$data[0,9] = "Ligand derived specificity of total expression value"
$data[0,10] = "Receptor-expressing cells"
$data[0,11] = "Receptor detection rate"
$data[0,12] = "Receptor average expression value"
$data[0,13] = "Receptor total expression value"
$data[0,14] = "Receptor derived specificity of average expression value"
$data[0,15] = "Receptor derived specificity of total expression value"
$data[0,16] = "Edge average expression weight"
$data[0,17] = "Edge total expression weight"
$data[0,18] = "Edge average expression derived specificity"
$data[0,19] = "Edge total expression derived specificity"
$data[1,0] = "ECs"
$data[1,1] = "Slit1"
$data[1,2] = "Robo2"
$data[1,3] = "ECs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.08454100000000002
$data[1,7] = 0.253623
$data[1,8] = 0.04188307112135965
$data[1,9] = 0.04356611262330446
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.3615393333333333
$data[1,13] = 1.084618
$data[1,14] = 0.7649240942154193
$data[1,15] = 0.7664955283791567
$data[1,16] = 0.03056489677933334
$data[1,17] = 0.275084071014
$data[1,18] = 0.03203737024046602
$data[1,19] = 0.0333932305146256
$data[2,0] = "ECs"
$data[2,1] = "Slit1"
$data[2,2] = "Robo2"
$data[2,3] = "FAPs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.08454100000000002
$data[2,7] = 0.253623
$data[2,8] = 0.04188307112135965
$data[2,9] = 0.04356611262330446
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.108201
$data[2,13] = 0.324603
$data[2,14] = 0.228925442648571
$data[2,15] = 0.2293957393280025
$data[2,16] = 0.009147420741
$data[2,17] = 0.08232678666900001
$data[2,18] = 0.009588100595938839
$data[2,19] = 0.009993880614869948
$data[3,0] = "ECs"
$data[3,1] = "Slit1"
$data[3,2] = "Robo2"
$data[3,3] = "MuSCs"
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.08454100000000002
$data[3,7] = 0.253623
$data[3,8] = 0.04188307112135965
$data[3,9] = 0.04356611262330446
$data[3,10] = 1
$data[3,11] = 0.5
$data[3,12] = 0.002907
$data[3,13] = 0.005814
$data[3,14] = 0.006150463136009796
$data[3,15] = 0.004108732292840814
$data[3,16] = 0.0002457606870000001
$data[3,17] = 0.001474564122
$data[3,18] = 0.000257600284954799
$data[3,19] = 0.0001790014938089108
$data[4,0] = "FAPs"
$data[4,1] = "Slit1"
$data[4,2] = "Robo2"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.340305
$data[4,7] = 1.020915
$data[4,8] = 0.1685929728528678
$data[4,9] = 0.1753677618702597
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.3615393333333333
$data[4,13] = 1.084618
$data[4,14] = 0.7649240942154193
$data[4,15] = 0.7664955283791567
$data[4,16] = 0.12303364283
$data[4,17] = 1.10730278547
$data[4,18] = 0.1289608270505646
$data[4,19] = 0.1344186052954148
$data[5,0] = "FAPs"
$data[5,1] = "Slit1"
$data[5,2] = "Robo2"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.340305
$data[5,7] = 1.020915
$data[5,8] = 0.1685929728528678
$data[5,9] = 0.1753677618702597
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.108201
$data[5,13] = 0.324603
$data[5,14] = 0.228925442648571
$data[5,15] = 0.2293957393280025
$data[5,16] = 0.036821341305
$data[5,17] = 0.331392071745
$data[5,18] = 0.03859522093778126
$data[5,19] = 0.0402286173885253
$data[6,0] = "FAPs"
$data[6,1] = "Slit1"
$data[6,2] = "Robo2"
$data[6,3] = "MuSCs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.340305
$data[6,7] = 1.020915
$data[6,8] = 0.1685929728528678
$data[6,9] = 0.1753677618702597
$data[6,10] = 1
$data[6,11] = 0.5
$data[6,12] = 0.002907
$data[6,13] = 0.005814
$data[6,14] = 0.006150463136009796
$data[6,15] = 0.004108732292840814
$data[6,16] = 0.000989266635
$data[6,17] = 0.00593559981
$data[6,18] = 0.001036924864521863
$data[6,19] = 0.0007205391863195537
$data[7,0] = "Inflammatory-Mac"
$data[7,1] = "Slit1"
$data[7,2] = "Robo2"
$data[7,3] = "ECs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.4435069999999999
$data[7,7] = 1.330521
$data[7,8] = 0.2197210255830999
$data[7,9] = 0.2285503591301722
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.3615393333333333
$data[7,13] = 1.084618
$data[7,14] = 0.7649240942154193
$data[7,15] = 0.7664955283791567
$data[7,16] = 0.1603452251086666
$data[7,17] = 1.443107025978
$data[7,18] = 0.1680699064742356
$data[7,19] = 0.1751828282827273
$data[8,0] = "Inflammatory-Mac"
$data[8,1] = "Slit1"
$data[8,2] = "Robo2"
$data[8,3] = "FAPs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.4435069999999999
$data[8,7] = 1.330521
$data[8,8] = 0.2197210255830999
$data[8,9] = 0.2285503591301722
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.108201
$data[8,13] = 0.324603
$data[8,14] = 0.228925442648571
$data[8,15] = 0.2293957393280025
$data[8,16] = 0.04798790090699999
$data[8,17] = 0.4318911081629999
$data[8,18] = 0.05029973304080913
$data[8,19] = 0.05242847860634633
$data[9,0] = "Inflammatory-Mac"
$data[9,1] = "Slit1"
$data[9,2] = "Robo2"
$data[9,3] = "MuSCs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.4435069999999999
$data[9,7] = 1.330521
$data[9,8] = 0.2197210255830999
$data[9,9] = 0.2285503591301722
$data[9,10] = 1
$data[9,11] = 0.5
$data[9,12] = 0.002907
$data[9,13] = 0.005814
$data[9,14] = 0.006150463136009796
$data[9,15] = 0.004108732292840814
$data[9,16] = 0.001289274849
$data[9,17] = 0.007735649093999999
$data[9,18] = 0.001351386068055121
$data[9,19] = 0.0009390522410985036
$data[10,0] = "MuSCs"
$data[10,1] = "Slit1"
$data[10,2] = "Robo2"
$data[10,3] = "ECs"
$data[10,4] = 2
$data[10,5] = 1
$data[10,6] = 0.2339355
$data[10,7] = 0.467871
$data[10,8] = 0.1158956859312148
$data[10,9] = 0.08036858123741962
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.3615393333333333
$data[10,13] = 1.084618
$data[10,14] = 0.7649240942154193
$data[10,15] = 0.7664955283791567
$data[10,16] = 0.08457688471299998
$data[10,17] = 0.507461308278
$data[10,18] = 0.08865140258440918
$data[10,19] = 0.06160215814065913
$data[11,0] = "MuSCs"
$data[11,1] = "Slit1"
$data[11,2] = "Robo2"
$data[11,3] = "FAPs"
$data[11,4] = 2
$data[11,5] = 1
$data[11,6] = 0.2339355
$data[11,7] = 0.467871
$data[11,8] = 0.1158956859312148
$data[11,9] = 0.08036858123741962
$data[11,10] = 2
$data[11,11] = 0.6666666666666666
$data[11,12] = 0.108201
$data[11,13] = 0.324603
$data[11,14] = 0.228925442648571
$data[11,15] = 0.2293957393280025
$data[11,16] = 0.0253120550355
$data[11,17] = 0.151872330213
$data[11,18] = 0.0265314712028631
$data[11,19] = 0.0184362101117005
$data[12,0] = "MuSCs"
$data[12,1] = "Slit1"
$data[12,2] = "Robo2"
$data[12,3] = "MuSCs"
$data[12,4] = 2
$data[12,5] = 1
$data[12,6] = 0.2339355
$data[12,7] = 0.467871
$data[12,8] = 0.1158956859312148
$data[12,9] = 0.08036858123741962
$data[12,10] = 1
$data[12,11] = 0.5
$data[12,12] = 0.002907
$data[12,13] = 0.005814
$data[12,14] = 0.006150463136009796
$data[12,15] = 0.004108732292840814
$data[12,16] = 0.0006800504984999999
$data[12,17] = 0.002720201994
$data[12,18] = 0.0007128121439425055
$data[12,19] = 0.0003302129850599863
$data[13,0] = "Resolving-Mac"
$data[13,1] = "Slit1"
$data[13,2] = "Robo2"
$data[13,3] = "ECs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 0.9162119999999999
$data[13,7] = 2.748636
$data[13,8] = 0.4539072445114579
$data[13,9] = 0.4721471851388441
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 0.3615393333333333
$data[13,13] = 1.084618
$data[13,14] = 0.7649240942154193
$data[13,15] = 0.7664955283791567
$data[13,16] = 0.3312466756719999
$data[13,17] = 2.981220081047999
$data[13,18] = 0.3472045878657438
$data[13,19] = 0.3618987061457298
$data[14,0] = "Resolving-Mac"
$data[14,1] = "Slit1"
$data[14,2] = "Robo2"
$data[14,3] = "FAPs"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 0.9162119999999999
$data[14,7] = 2.748636
$data[14,8] = 0.4539072445114579
$data[14,9] = 0.4721471851388441
$data[14,10] = 2
$data[14,11] = 0.6666666666666666
$data[14,12] = 0.108201
$data[14,13] = 0.324603
$data[14,14] = 0.228925442648571
$data[14,15] = 0.2293957393280025
$data[14,16] = 0.09913505461199998
$data[14,17] = 0.8922154915079998
$data[14,18] = 0.1039109168711786
$data[14,19] = 0.1083085526065604
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Slit1"
$data[15,2] = "Robo2"
$data[15,3] = "MuSCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 0.9162119999999999
$data[15,7] = 2.748636
$data[15,8] = 0.4539072445114579
$data[15,9] = 0.4721471851388441
$data[15,10] = 1
$data[15,11] = 0.5
$data[15,12] = 0.002907
$data[15,13] = 0.005814
$data[15,14] = 0.006150463136009796
$data[15,15] = 0.004108732292840814
$data[15,16] = 0.002663428284
$data[15,17] = 0.015980569704
$data[15,18] = 0.002791739774535506
$data[15,19] = 0.001939926386553859

$ws.Range("A1:T16").Value = $data
